$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.406.16"
$ws.Cells.Item(2, 5).Value = "  +2.01%  "
$ws.Cells.Item(3, 4).Value = "1.844.65"
$ws.Cells.Item(3, 5).Value = "  +1.71%  "
$ws.Cells.Item(4, 4).Value = "'1.014"
$ws.Cells.Item(4, 5).Value = "  +1.32%  "
$ws.Cells.Item(5, 4).Value = "'315.88"
$ws.Cells.Item(5, 5).Value = "  +2.33%  "
$ws.Cells.Item(6, 4).Value = "'1.012"
$ws.Cells.Item(7, 4).Value = "'0.4744"
$ws.Cells.Item(7, 5).Value = "  +1.67%  "
$ws.Cells.Item(8, 4).Value = "'0.3704"
$ws.Cells.Item(8, 5).Value = "  +0.57%  "
$ws.Cells.Item(9, 4).Value = "'0.07475"
$ws.Cells.Item(9, 5).Value = "  +1.44%  "
$ws.Cells.Item(10, 4).Value = "'0.8879"
$ws.Cells.Item(10, 5).Value = "  +2.13%  "
$ws.Cells.Item(11, 4).Value = "'20.54"
$ws.Cells.Item(11, 5).Value = "  +0.77%  "
$ws.Cells.Item(12, 4).Value = "1.857.63"
$ws.Cells.Item(12, 5).Value = "  +2.40%  "
$ws.Cells.Item(13, 4).Value = "'0.07401"
$ws.Cells.Item(13, 5).Value = "  +4.50%  "
$ws.Cells.Item(14, 4).Value = "'5.494"
$ws.Cells.Item(14, 5).Value = "  +2.69%  "
$ws.Cells.Item(15, 4).Value = "'93.33"
$ws.Cells.Item(15, 5).Value = "  +1.83%  "
$ws.Cells.Item(16, 5).Value = "  +1.49%  "
$ws.Cells.Item(17, 4).Value = "'1.014"
$ws.Cells.Item(17, 5).Value = "  +1.31%  "
$ws.Cells.Item(18, 4).Value = "'0.000008871"
$ws.Cells.Item(18, 5).Value = "  +2.08%  "
$ws.Cells.Item(19, 4).Value = "'1.013"
$ws.Cells.Item(19, 5).Value = "  +1.16%  "
$ws.Cells.Item(20, 4).Value = "'14.88"
$ws.Cells.Item(20, 5).Value = "  +1.03%  "
$ws.Cells.Item(21, 4).Value = "27.429.64"
$ws.Cells.Item(21, 5).Value = "  +1.93%  "
$ws.Cells.Item(22, 4).Value = "'5.352"
$ws.Cells.Item(22, 5).Value = "  +0.28%  "
$ws.Cells.Item(23, 5).Value = "  +1.63%  "
$ws.Cells.Item(24, 4).Value = "2.077.09"
$ws.Cells.Item(24, 5).Value = "  +1.12%  "
$ws.Cells.Item(25, 5).Value = "  +0.56%  "
$ws.Cells.Item(26, 5).Value = "  +0.92%  "
$ws.Cells.Item(27, 5).Value = "  +1.57%  "
$ws.Cells.Item(28, 4).Value = "'2.184"
$ws.Cells.Item(28, 5).Value = "  +0.67%  "
$ws.Cells.Item(29, 4).Value = "'5.288"
$ws.Cells.Item(29, 5).Value = "  -0.52%  "
$ws.Cells.Item(30, 4).Value = "'118.15"
$ws.Cells.Item(30, 5).Value = "  +2.33%  "
$ws.Cells.Item(31, 4).Value = "'0.08979"
$ws.Cells.Item(31, 5).Value = "  +0.52%  "
$ws.Cells.Item(32, 4).Value = "'0.7620"
$ws.Cells.Item(32, 5).Value = "  -0.58%  "
$ws.Cells.Item(33, 5).Value = "  +1.70%  "
$ws.Cells.Item(34, 4).Value = "'4.569"
$ws.Cells.Item(34, 5).Value = "  +1.54%  "
$ws.Cells.Item(35, 4).Value = "'2.948"
$ws.Cells.Item(35, 5).Value = "  +1.60%  "
$ws.Cells.Item(36, 4).Value = "'1.013"
$ws.Cells.Item(37, 4).Value = "'1.109"
$ws.Cells.Item(37, 5).Value = "  +1.99%  "
$ws.Cells.Item(38, 4).Value = "'0.05370"
$ws.Cells.Item(38, 5).Value = "  +1.71%  "
$ws.Cells.Item(39, 4).Value = "'0.01968"
$ws.Cells.Item(39, 5).Value = "  +0.42%  "
$ws.Cells.Item(40, 4).Value = "'3.006"
$ws.Cells.Item(40, 5).Value = "  +2.41%  "
$ws.Cells.Item(41, 4).Value = "'7.330"
$ws.Cells.Item(41, 5).Value = "  +0.97%  "
$ws.Cells.Item(42, 4).Value = "'2.396"
$ws.Cells.Item(42, 5).Value = "  +1.58%  "
$ws.Cells.Item(43, 4).Value = "'0.5361"
$ws.Cells.Item(43, 5).Value = "  +0.82%  "
$ws.Cells.Item(44, 5).Value = "  +0.26%  "
$ws.Cells.Item(45, 4).Value = "'8.563"
$ws.Cells.Item(45, 5).Value = "  +1.79%  "
$ws.Cells.Item(46, 4).Value = "'0.4979"
$ws.Cells.Item(46, 5).Value = "  +1.14%  "
$ws.Cells.Item(47, 4).Value = "'10.57"
$ws.Cells.Item(47, 5).Value = "  +0.91%  "
$ws.Cells.Item(49, 4).Value = "'105.24"
$ws.Cells.Item(49, 5).Value = "  +1.51%  "
$ws.Cells.Item(50, 4).Value = "'1.684"
$ws.Cells.Item(50, 5).Value = "  +1.10%  "
$ws.Cells.Item(51, 4).Value = "'0.06331"
$ws.Cells.Item(51, 5).Value = "  +0.69%  "
